{"js": "// The \"Requisitos\" list paragraph lists three prerequisite courses, each\n// rendered as its own run (text + manual line break <w:br/>). The edit\n// moves the \"LOB1019 -  F\u00edsica II  (Requisito fraco)\" line from the first\n// position to the last position, leaving the other two lines' relative\n// order unchanged:\n//   Before: LOB1019, LOQ4053, LOB1004\n//   After:  LOQ4053, LOB1004, LOB1019\n\nconst LINE_TO_MOVE = \"LOB1019 -  F\u00edsica II  (Requisito fraco)\";\nconst ANCHOR_TEXT = \"Requisito fraco\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the requirements paragraph (the ListBullet paragraph that holds the\n// three \"Requisito fraco\" lines).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(ANCHOR_TEXT) >= 0 && p.text.indexOf(LINE_TO_MOVE) >= 0) {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the Requisitos paragraph containing '\" + LINE_TO_MOVE + \"'\");\n}\n\n// Search within that paragraph for the line's text INCLUDING its trailing\n// manual line break (\\v) so the whole run (text + <w:br/>) is matched.\nconst hits = target.search(LINE_TO_MOVE + \"\\u000b\", { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the line to move inside the Requisitos paragraph\");\n}\n\nconst lineRange = hits.items[0];\nconst lineTextWithBreak = lineRange.text; // \"LOB1019 -  F\u00edsica II  (Requisito fraco)\\u000b\"\n\n// Remove the line (and its break) from its current position...\nlineRange.delete();\nawait context.sync();\n\n// ...and re-append it at the end of the paragraph so it becomes the last\n// line, landing in its own new run exactly like the other two.\ntarget.insertText(lineTextWithBreak, \"End\");\nawait context.sync();\n", "ps1": "# The \"Requisitos\" list paragraph lists three prerequisite courses, each\n# rendered as its own run (text + manual line break vertical-tab mark).\n# The edit moves the \"LOB1019 -  F\u00edsica II  (Requisito fraco)\" line from\n# the first position to the last position, leaving the other two lines'\n# relative order unchanged:\n#   Before: LOB1019, LOQ4053, LOB1004\n#   After:  LOQ4053, LOB1004, LOB1019\n\n$d = $word.ActiveDocument\n\n$lineToMove = \"LOB1019 -  F\u00edsica II  (Requisito fraco)\"\n\n# Find the requirements paragraph: the one containing all three\n# \"Requisito fraco\" lines, including the line we need to move.\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -match \"Requisito fraco\" -and $t.Contains($lineToMove)) {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the Requisitos paragraph containing '$lineToMove'\"\n}\n\n$r = $target.Range\n$pStart = $r.Start\n$pEnd = $r.End\n$fullText = $r.Text\n\n# Locate the line inside the paragraph text; the break character\n# (manual line break, chr 11) immediately follows it.\n$relIdx = $fullText.IndexOf($lineToMove)\nif ($relIdx -lt 0) {\n    throw \"Could not locate '$lineToMove' inside the Requisitos paragraph\"\n}\n$lineLenWithBreak = $lineToMove.Length + 1\n\n$lineStart = $pStart + $relIdx\n$lineEnd = $lineStart + $lineLenWithBreak\n\n$lineRange = $d.Range($lineStart, $lineEnd)\n$lineText = $lineRange.Text\n\n# Remove the line (and its trailing break) from its current position.\n$lineRange.Delete()\n\n# The paragraph end moved back by the removed length; compute the new\n# insertion point, which sits right before the paragraph mark.\n$newParaEnd = $pEnd - $lineLenWithBreak\n$insertionPoint = $d.Range($newParaEnd - 1, $newParaEnd - 1)\n$insertionPoint.InsertAfter($lineText)\n\n$d.Save()\n"}
